$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = [double]"0.2630136783801007"
$ws.Range("E2").Value = [double]"0.2630136783801007"

# Row 3
$ws.Range("D3").Value = [double]"0.9999969903038608"
$ws.Range("E3").Value = [double]"0.9999969903038608"

# Row 4
$ws.Range("D4").Value = [double]"6.247737606847879E-26"
$ws.Range("E4").Value = [double]"6.247737606847879E-26"

# Row 5
$ws.Range("D5").Value = [double]"0.9999998725965373"
$ws.Range("E5").Value = [double]"0.9999998725965373"

# Row 6
$ws.Range("D6").Value = [double]"0.996379558015915"
$ws.Range("E6").Value = [double]"0.996379558015915"

# Row 7
$ws.Range("D7").Value = [double]"7.876848809716225E-08"
$ws.Range("E7").Value = [double]"0.9999999212315119"

# Row 9
$ws.Range("D9").Value = [double]"0.9999995583469222"
$ws.Range("E9").Value = [double]"4.416530777762162E-07"

# Row 10
$ws.Range("D10").Value = [double]"0.02316974730604013"
$ws.Range("E10").Value = [double]"0.9768302526939598"

# Row 11
$ws.Range("D11").Value = [double]"2.230715265620821E-06"
$ws.Range("E11").Value = [double]"0.9999977692847344"
$ws.Range("F11").Value = [double]"6.765077114105225"
$ws.Range("G11").Value = [double]"0.4"
